# Build v2.1.2: Fix SearchCriteria variants and Schemas sheet grouping/sorting
#
# - "Body" sheet: row 3 (the request body's "dateTime" field) is replaced by a
#   single "schema" reference to the generated request schema
#   (retransmitOutFiles.210702Request); the old "networkFileName" row (row 4)
#   is removed.
# - "200" sheet: row 3 (the response's "dateTime" field) is replaced by a
#   single "schema" reference to the generated response schema
#   (retransmitOutFiles.210702Response); the old "commandRef"/"commandStatus"
#   rows (4 and 5) are removed.
# - "204" sheet: gains a new row 3 that is a "schema" reference to the same
#   response schema (retransmitOutFiles.210702Response).
# - "400" sheet: row 3 (the response's "dateTime" field) is replaced by a
#   single "schema" reference to a shared "errorResponse" schema; the old
#   "errorCode"/"errorCodeDescription"/"requestId" rows (4,5,6) are removed.
# - "401", "403", "404", "429", "500" sheets: each gains a new row 3 that is a
#   "schema" reference to a shared "errorResponse1" schema.

$wb = $excel.ActiveWorkbook

function Set-SchemaRow {
    param(
        $ws,
        [int]$row,
        [string]$section,
        [string]$schemaName
    )

    $ws.Cells.Item($row, 1).Value = $section      # A - Section
    $ws.Cells.Item($row, 2).Value = $schemaName   # B - Name
    $ws.Cells.Item($row, 3).ClearContents()       # C - Parent
    $ws.Cells.Item($row, 4).ClearContents()       # D - Description
    $ws.Cells.Item($row, 5).Value = "schema"      # E - Type
    $ws.Cells.Item($row, 6).ClearContents()       # F - Items Data Type
    $ws.Cells.Item($row, 7).Value = $schemaName   # G - Schema Name
    $ws.Cells.Item($row, 8).ClearContents()       # H - Format
    $ws.Cells.Item($row, 9).Value = "Yes"         # I - Mandatory
    $ws.Cells.Item($row, 10).ClearContents()      # J - Min
    $ws.Cells.Item($row, 11).ClearContents()      # K - Max
    $ws.Cells.Item($row, 12).ClearContents()      # L - PatternEba
    $ws.Cells.Item($row, 13).ClearContents()      # M - Regex
    $ws.Cells.Item($row, 14).ClearContents()      # N - Allowed value
    $ws.Cells.Item($row, 15).ClearContents()      # O - Example
}

# --- Body: request schema, drop the old networkFileName row -----------------
$wsBody = $wb.Worksheets.Item("Body")
Set-SchemaRow $wsBody 3 "body" "retransmitOutFiles.210702Request"
$wsBody.Range("A4:A4").EntireRow.Delete()

# --- 200: response schema, drop the old commandRef/commandStatus rows -------
$ws200 = $wb.Worksheets.Item("200")
Set-SchemaRow $ws200 3 "content" "retransmitOutFiles.210702Response"
$ws200.Range("A4:A5").EntireRow.Delete()

# --- 204: add a new response-schema row -------------------------------------
$ws204 = $wb.Worksheets.Item("204")
Set-SchemaRow $ws204 3 "content" "retransmitOutFiles.210702Response"

# --- 400: shared errorResponse schema, drop errorCode/.../requestId rows ----
$ws400 = $wb.Worksheets.Item("400")
Set-SchemaRow $ws400 3 "content" "errorResponse"
$ws400.Range("A4:A6").EntireRow.Delete()

# --- 401, 403, 404, 429, 500: add a new errorResponse1-schema row ----------
foreach ($sheetName in @("401", "403", "404", "429", "500")) {
    $ws = $wb.Worksheets.Item($sheetName)
    Set-SchemaRow $ws 3 "content" "errorResponse1"
}
